# Sample Project / Main.xlsx - "Rules" sheet: B11 used to hold the shared
# string "R40" (last greeting-range rule). The commit replaces that cell's
# content with the text "1" (kept as a genuine text value, not a number),
# while leaving every other cell, and B11's own style (s="23"), untouched.
#
# Excel auto-converts a literal digit string typed/assigned straight into a
# General-formatted cell into a number (losing the shared-string "t=s"
# typing the diff expects). To reproduce a real text "1" without disturbing
# B11's existing number format/style, we stage the text in a scratch cell
# that we temporarily mark as Text, copy *values* (not formats) from there
# onto B11, then restore the scratch cell to its original blank state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("B5")     # already-blank cell inside the used range
$parking = $ws.Range("H20")    # far outside the used range - safe holding spot

# Remember the scratch cell's original formatting so it can be restored.
$scratch.Copy()
$parking.PasteSpecial(-4122)   # xlPasteFormats

# Stage the literal text "1" in the scratch cell (Text format forces Excel
# to keep it as a string instead of coercing it to a number).
$scratch.NumberFormat = "@"
$scratch.Value = "1"

# Copy only the value into B11 so B11 keeps its own existing style (s="23").
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues

# Put the scratch cell back exactly the way it was (blank, original format).
$scratch.ClearContents()
$parking.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats
$parking.Clear()
